# Generate Report for Handoff
# Marks the six "Ready for handoff" rows (the e2e files whose status is
# "Ready for handoff") as handed-off with priority "ht" and refreshes the
# handoff timestamps on all three report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Rows 7,8,9,10,11,13 correspond to the files that were just handed off.
$rows = @(7, 8, 9, 10, 11, 13)

foreach ($r in $rows) {
    # Priority column (E) gets marked as handoff type "ht".
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Latest Handoff Datetime (column H) is refreshed per locale.
    $zhcn.Range("H$r").Value = "2016-08-27 18:28:56"
    $dede.Range("H$r").Value = "2016-08-27 18:29:04"

    # Overview sheet's Latest HO Xliff Generate Date (column G) mirrors the
    # de-de handoff timestamp (the later of the two locale refreshes).
    $overview.Range("G$r").Value = "2016-08-27 18:29:04"
}
